$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped price/volume updates for this run. Every source value
# is plain scraped text (prices use "." as a thousands separator in this feed,
# so e.g. "1.632.70" is text, not a number) and must stay a text cell exactly
# like the rest of the column. Plain-decimal-looking prices (e.g. "213.88",
# "15.20") would otherwise be auto-coerced to numbers by Excel, so those are
# forced to keep a text format, then restored to the sheet's normal style.
$ws.Range("D2").Value = "25.983.43"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.632.70"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.859.80"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "1.665.61"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₃0745"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "25.986.49"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "1.132.78"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.865"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("D44").Value = "1.770.34"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("E51").Value = "  +0.32%  "
